# Update NATMI LR-pair (L1cam-Itga5) TPM-derived metrics with recomputed values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.741029
$ws.Range("H2").Value = 23.223087
$ws.Range("I2").Value = 0.4930486933812723
$ws.Range("J2").Value = 0.4930486933812723
$ws.Range("M2").Value = 11.146846
$ws.Range("N2").Value = 33.440538
$ws.Range("O2").Value = 0.2594806085672136
$ws.Range("P2").Value = 0.2594806085672136
$ws.Range("Q2").Value = 86.28805814453402
$ws.Range("R2").Value = 776.5925233008061
$ws.Range("S2").Value = 0.127936575011842
$ws.Range("T2").Value = 0.127936575011842
$ws.Range("G3").Value = 7.741029
$ws.Range("H3").Value = 23.223087
$ws.Range("I3").Value = 0.4930486933812723
$ws.Range("J3").Value = 0.4930486933812723
$ws.Range("O3").Value = 0.6444737471070977
$ws.Range("P3").Value = 0.6444737471070977
$ws.Range("Q3").Value = 214.314235156413
$ws.Range("R3").Value = 1928.828116407717
$ws.Range("S3").Value = 0.317756938929687
$ws.Range("T3").Value = 0.317756938929687
$ws.Range("G4").Value = 7.741029
$ws.Range("H4").Value = 23.223087
$ws.Range("I4").Value = 0.4930486933812723
$ws.Range("J4").Value = 0.4930486933812723
$ws.Range("O4").Value = 0.09604564432568881
$ws.Range("P4").Value = 0.09604564432568881
$ws.Range("Q4").Value = 31.93915795043901
$ws.Range("R4").Value = 287.452421553951
$ws.Range("S4").Value = 0.04735517943974327
$ws.Range("T4").Value = 0.04735517943974327
$ws.Range("I5").Value = 0.0194007766416684
$ws.Range("J5").Value = 0.0194007766416684
$ws.Range("M5").Value = 11.146846
$ws.Range("N5").Value = 33.440538
$ws.Range("O5").Value = 0.2594806085672136
$ws.Range("P5").Value = 0.2594806085672136
$ws.Range("Q5").Value = 3.395314429138667
$ws.Range("R5").Value = 30.55782986224801
$ws.Range("S5").Value = 0.005034125329656698
$ws.Range("T5").Value = 0.005034125329656698
$ws.Range("I6").Value = 0.0194007766416684
$ws.Range("J6").Value = 0.0194007766416684
$ws.Range("O6").Value = 0.6444737471070977
$ws.Range("P6").Value = 0.6444737471070977
$ws.Range("S6").Value = 0.01250329121904389
$ws.Range("T6").Value = 0.01250329121904389
$ws.Range("I7").Value = 0.0194007766416684
$ws.Range("J7").Value = 0.0194007766416684
$ws.Range("O7").Value = 0.09604564432568881
$ws.Range("P7").Value = 0.09604564432568881
$ws.Range("S7").Value = 0.001863360092967815
$ws.Range("T7").Value = 0.001863360092967815
$ws.Range("G8").Value = 7.654706000000001
$ws.Range("I8").Value = 0.4875505299770593
$ws.Range("J8").Value = 0.4875505299770593
$ws.Range("M8").Value = 11.146846
$ws.Range("N8").Value = 33.440538
$ws.Range("O8").Value = 0.2594806085672136
$ws.Range("P8").Value = 0.2594806085672136
$ws.Range("Q8").Value = 85.32582895727603
$ws.Range("R8").Value = 767.9324606154842
$ws.Range("S8").Value = 0.1265099082257148
$ws.Range("T8").Value = 0.1265099082257148
$ws.Range("G9").Value = 7.654706000000001
$ws.Range("I9").Value = 0.4875505299770593
$ws.Range("J9").Value = 0.4875505299770593
$ws.Range("O9").Value = 0.6444737471070977
$ws.Range("P9").Value = 0.6444737471070977
$ws.Range("S9").Value = 0.3142135169583668
$ws.Range("T9").Value = 0.3142135169583668
$ws.Range("G10").Value = 7.654706000000001
$ws.Range("I10").Value = 0.4875505299770593
$ws.Range("J10").Value = 0.4875505299770593
$ws.Range("O10").Value = 0.09604564432568881
$ws.Range("P10").Value = 0.09604564432568881
$ws.Range("S10").Value = 0.04682710479297772
$ws.Range("T10").Value = 0.04682710479297772

Write-Output "TPM values updated"
